$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D and E to be treated as plain text so that numeric-looking
# strings (e.g. "239.05", "65.60") are stored verbatim instead of being
# coerced into floating point numbers by Excel's automatic type detection.
$dataRange = $ws.Range("D2:E51")
$dataRange.NumberFormat = "@"

$ws.Range("D2").Value = "31.080.76"
$ws.Range("E2").Value = "  +2.86%  "

$ws.Range("D3").Value = "1.898.11"
$ws.Range("E3").Value = "  +3.18%  "

$ws.Range("D4").Value = "0.9973"
$ws.Range("E4").Value = "  -0.25%  "

$ws.Range("D5").Value = "239.05"
$ws.Range("E5").Value = "  +2.77%  "

$ws.Range("D6").Value = "0.9974"
$ws.Range("E6").Value = "  -0.26%  "

$ws.Range("E7").Value = "  +2.66%  "

$ws.Range("D8").Value = "0.2862"
$ws.Range("E8").Value = "  +5.56%  "

$ws.Range("D9").Value = "0.06551"
$ws.Range("E9").Value = "  +4.41%  "

$ws.Range("D10").Value = "18.88"
$ws.Range("E10").Value = "  +17.30%  "

$ws.Range("D11").Value = "1.868.24"
$ws.Range("E11").Value = "  +1.55%  "

$ws.Range("D12").Value = "96.10"
$ws.Range("E12").Value = "  +14.77%  "

$ws.Range("D13").Value = "0.07525"
$ws.Range("E13").Value = "  +1.50%  "

$ws.Range("D14").Value = "5.152"
$ws.Range("E14").Value = "  +4.54%  "

$ws.Range("D15").Value = "0.6549"
$ws.Range("E15").Value = "  +5.55%  "

$ws.Range("D16").Value = "298.41"
$ws.Range("E16").Value = "  +31.58%  "

$ws.Range("D17").Value = "31.065.33"
$ws.Range("E17").Value = "  +3.09%  "

$ws.Range("D18").Value = "13.20"
$ws.Range("E18").Value = "  +6.77%  "

$ws.Range("D19").Value = "0.9970"
$ws.Range("E19").Value = "  -0.28%  "

$ws.Range("D20").Value = "0.000007588"
$ws.Range("E20").Value = "  +4.04%  "

$ws.Range("D21").Value = "0.9966"
$ws.Range("E21").Value = "  -0.35%  "

$ws.Range("D22").Value = "5.215"
$ws.Range("E22").Value = "  +6.78%  "

$ws.Range("E23").Value = "  +5.50%  "

$ws.Range("D24").Value = "9.335"
$ws.Range("E24").Value = "  +1.51%  "

$ws.Range("D25").Value = "168.47"
$ws.Range("E25").Value = "  +2.53%  "

$ws.Range("E26").Value = "  +10.62%  "

$ws.Range("E27").Value = "  +5.37%  "

$ws.Range("D28").Value = "0.1061"
$ws.Range("E28").Value = "  +1.82%  "

$ws.Range("D29").Value = "1.364"
$ws.Range("E29").Value = "  -0.45%  "

$ws.Range("D30").Value = "4.171"
$ws.Range("E30").Value = "  +2.31%  "

$ws.Range("D31").Value = "3.994"
$ws.Range("E31").Value = "  +5.30%  "

$ws.Range("D32").Value = "0.05022"
$ws.Range("E32").Value = "  +4.30%  "

$ws.Range("D33").Value = "1.189"
$ws.Range("E33").Value = "  +4.35%  "

$ws.Range("D34").Value = "0.7289"
$ws.Range("E34").Value = "  +2.82%  "

$ws.Range("D35").Value = "2.709"
$ws.Range("E35").Value = "  +0.53%  "

$ws.Range("D36").Value = "0.01944"
$ws.Range("E36").Value = "  +3.72%  "

$ws.Range("D37").Value = "2.724"
$ws.Range("E37").Value = "  +2.76%  "

$ws.Range("D38").Value = "2.069"
$ws.Range("E38").Value = "  +7.77%  "

$ws.Range("D39").Value = "0.9028"
$ws.Range("E39").Value = "  +1.24%  "

$ws.Range("D40").Value = "107.68"
$ws.Range("E40").Value = "  +3.32%  "

$ws.Range("D41").Value = "0.4246"
$ws.Range("E41").Value = "  +6.00%  "

$ws.Range("D42").Value = "0.9962"
$ws.Range("E42").Value = "  -0.53%  "

$ws.Range("D43").Value = "5.606"
$ws.Range("E43").Value = "  +1.35%  "

$ws.Range("D44").Value = "7.418"
$ws.Range("E44").Value = "  +5.36%  "

$ws.Range("D45").Value = "65.60"
$ws.Range("E45").Value = "  +9.69%  "

$ws.Range("D46").Value = "0.1232"
$ws.Range("E46").Value = "  +3.13%  "

$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").Value = "8.953"
$ws.Range("E47").Value = "  +5.23%  "

$ws.Range("B48").Value = "Elrond"
$ws.Range("C48").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D48").Value = "34.76"
$ws.Range("E48").Value = "  +6.16%  "

$ws.Range("D49").Value = "1.403"
$ws.Range("E49").Value = "  +3.55%  "

$ws.Range("D50").Value = "0.05592"
$ws.Range("E50").Value = "  +1.51%  "

$ws.Range("D51").Value = "0.3815"
$ws.Range("E51").Value = "  +4.93%  "

# Restore the default ("Normal") style on the data range so no residual
# number-format / quote-prefix style is left on the cells.
$dataRange.Style = "Normal"